$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Status" row -> flip from active to draft
$ws.Range("A6").Value = "Status"
$ws.Range("B6").Value = "draft"

# "Date" row -> bump the publication timestamp
$ws.Range("A8").Value = "Date"
$ws.Range("B8").Value = "2023-08-01T16:12:28+00:00"

# Header row and body rows keep borders/fill/font; make sure the
# (already-present) top/wrap alignment is flagged as applied.
$ws.Range("A1:B1").Style.IncludeAlignment = $true
$ws.Range("A2:B21").Style.IncludeAlignment = $true
